$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.62%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.04%"

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.21%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07850"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.16%"

# Row 6
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.899"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.65%"

# Row 7
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.230"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.27%"

# Row 8
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.000"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.28%"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9379"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.43%"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1115"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-7.30%"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1946"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.42%"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09014"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.19%"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03344"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.36%"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09591"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.18%"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001378"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.57%"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006159"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.48%"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.608"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.61%"

# Row 18
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.431"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.99%"

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.11%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.424"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "27.65%"

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.14%"

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-10.52%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04390"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.29%"

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.80%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004567"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.58%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02223"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.45%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05042"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.95%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007448"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.74%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1350"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.66%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008740"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.53%"

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.53%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008192"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-14.52%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006547"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.36%"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-40.72%"
